# Insert a new data row at row 146 (pushing existing rows 146:234 down to
# 147:235) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(146).Insert()

$ws.Range("A146").Value = 7
$ws.Range("B146").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C146").Value = "Ñuble"
$ws.Range("D146").Value = 44830
$ws.Range("D146").NumberFormat = $ws.Range("D147").NumberFormat
$ws.Range("E146").Value = 16
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100104
$ws.Range("H146").Value = "Frutos de pepita"
$ws.Range("I146").Value = 100104005
$ws.Range("J146").Value = "Pera"
$ws.Range("K146").Value = "Packham's Triumph"
$ws.Range("L146").Value = "Primera"
$ws.Range("M146").Value = 120
$ws.Range("N146").Value = 10000
$ws.Range("O146").Value = 11000
$ws.Range("P146").Value = 10500
$ws.Range("Q146").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R146").Value = "Provincia de Curicó"
$ws.Range("S146").Value = 656
$ws.Range("T146").Value = 16
